$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix cells after simulating additional games.
# Every other matrix cell remains 0 (no observed transition).

# Row 2
$ws.Range("C2").Value = 1

# Row 3
$ws.Range("C3").Value = 0.09090909090909091
$ws.Range("J3").Value = 0.1818181818181818
$ws.Range("P3").Value = 0.5454545454545454
$ws.Range("S3").Value = 0.1818181818181818

# Row 4
$ws.Range("P4").Value = 1

# Row 6
$ws.Range("F6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.2142857142857143
$ws.Range("O6").Value = 0.07142857142857142
$ws.Range("Q6").Value = 0.07142857142857142
$ws.Range("R6").Value = 0.1428571428571428
$ws.Range("S6").Value = 0.4285714285714285

# Row 7
$ws.Range("Q7").Value = 0.25
$ws.Range("S7").Value = 0.75

# Row 8
$ws.Range("F8").Value = 0.04545454545454546
$ws.Range("J8").Value = 0.1818181818181818
$ws.Range("Q8").Value = 0.4090909090909091
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.2727272727272727

# Row 9
$ws.Range("B9").Value = 0.1
$ws.Range("J9").Value = 0.15
$ws.Range("Q9").Value = 0.35
$ws.Range("R9").Value = 0.25
$ws.Range("S9").Value = 0.15

# Row 10
$ws.Range("B10").Value = 0.05405405405405406
$ws.Range("D10").Value = 0.02027027027027027
$ws.Range("F10").Value = 0.06081081081081081
$ws.Range("J10").Value = 0.1689189189189189
$ws.Range("O10").Value = 0.02702702702702703
$ws.Range("Q10").Value = 0.3783783783783784
$ws.Range("R10").Value = 0.06081081081081081
$ws.Range("S10").Value = 0.2297297297297297

# Row 11
$ws.Range("J11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.7142857142857143

# Row 12
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("L12").Value = 0.1428571428571428

# Row 13
$ws.Range("J13").Value = 1

# Row 15
$ws.Range("F15").Value = 0.04166666666666666
$ws.Range("H15").Value = 0.04166666666666666
$ws.Range("I15").Value = 0.25
$ws.Range("J15").Value = 0.375
$ws.Range("O15").Value = 0.08333333333333333
$ws.Range("S15").Value = 0.2083333333333333

# Row 16
$ws.Range("F16").Value = 0.1111111111111111
$ws.Range("H16").Value = 0.3333333333333333
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.2222222222222222
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.1111111111111111

# Row 17
$ws.Range("H17").Value = 0.0958904109589041
$ws.Range("I17").Value = 0.0684931506849315
$ws.Range("J17").Value = 0.6027397260273972
$ws.Range("K17").Value = 0.0410958904109589
$ws.Range("M17").Value = 0.0136986301369863
$ws.Range("O17").Value = 0.0273972602739726
$ws.Range("S17").Value = 0.1506849315068493

# Row 18
$ws.Range("I18").Value = 0.1666666666666667
$ws.Range("J18").Value = 0.5555555555555556
$ws.Range("K18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.2222222222222222

# Row 19
$ws.Range("H19").Value = 0.1392405063291139
$ws.Range("I19").Value = 0.0759493670886076
$ws.Range("J19").Value = 0.5316455696202531
$ws.Range("K19").Value = 0.01265822784810127
$ws.Range("M19").Value = 0.01265822784810127
$ws.Range("O19").Value = 0.1392405063291139
$ws.Range("S19").Value = 0.08860759493670886
